$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H28").Value = 639.1905
$ws_ALC.Range("I28").Value = 658.0526
$ws_ALC.Range("K28").Value = 658.0526
$ws_ALC.Range("M28").Value = -173.0526
$ws_ALC.Range("H41").Value = 1722.2858
$ws_ALC.Range("J41").Value = 1936.25
$ws_ALC.Range("L41").Value = 1936.25
$ws_ALC.Range("N41").Value = -2816.25
$ws_ALC.Range("H53").Value = 426.3684
$ws_ALC.Range("I53").Value = 405.82352
$ws_ALC.Range("J53").Value = 601
$ws_ALC.Range("K53").Value = 405.82352
$ws_ALC.Range("L53").Value = 601
$ws_ALC.Range("M53").Value = 231.17648
$ws_ALC.Range("N53").Value = -1875
$ws_ALC.Range("H107").Value = 1127.1724
$ws_ALC.Range("J107").Value = 3799.6667
$ws_ALC.Range("L107").Value = 3799.6667
$ws_ALC.Range("N107").Value = -7639.6667
$ws_ALC.Range("H112").Value = 3354.75
$ws_ALC.Range("J112").Value = 3823.1667
$ws_ALC.Range("L112").Value = 11469.5001
$ws_ALC.Range("N112").Value = -13685.5001
$ws_ALC.Range("H132").Value = 3871.6765
$ws_ALC.Range("J132").Value = 6991.5
$ws_ALC.Range("L132").Value = 20974.5
$ws_ALC.Range("N132").Value = -26034.5
$ws_ALC.Range("H141").Value = 4600.3887
$ws_ALC.Range("I141").Value = 4600.3887
$ws_ALC.Range("K141").Value = 13801.1661
$ws_ALC.Range("M141").Value = -8621.166100000002

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 2093.1667
$ws_ARM.Range("I2").Value = 1347
$ws_ARM.Range("K2").Value = 1347
$ws_ARM.Range("M2").Value = -1234
$ws_ARM.Range("H32").Value = 23110.926
$ws_ARM.Range("I32").Value = 18559.84
$ws_ARM.Range("K32").Value = 18559.84
$ws_ARM.Range("M32").Value = -18272.84
$ws_ARM.Range("H110").Value = 4057.1428
$ws_ARM.Range("I110").Value = 900
$ws_ARM.Range("J110").Value = 4583.3335
$ws_ARM.Range("K110").Value = 900
$ws_ARM.Range("L110").Value = 4583.3335
$ws_ARM.Range("M110").Value = 1145
$ws_ARM.Range("N110").Value = -8673.333500000001
$ws_ARM.Range("H116").Value = 2093.1667
$ws_ARM.Range("I116").Value = 1347
$ws_ARM.Range("K116").Value = 1347
$ws_ARM.Range("M116").Value = 947

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 2093.1667
$ws_BSM.Range("I3").Value = 1347
$ws_BSM.Range("K3").Value = 1347
$ws_BSM.Range("M3").Value = -1233
$ws_BSM.Range("H29").Value = 1014.5
$ws_BSM.Range("I29").Value = 1014.5
$ws_BSM.Range("K29").Value = 1014.5
$ws_BSM.Range("M29").Value = -725.5
$ws_BSM.Range("H76").Value = 15813.5
$ws_BSM.Range("J76").Value = 15813.5
$ws_BSM.Range("L76").Value = 15813.5
$ws_BSM.Range("N76").Value = -16443.5
$ws_BSM.Range("H79").Value = 15813.5
$ws_BSM.Range("J79").Value = 15813.5
$ws_BSM.Range("L79").Value = 15813.5
$ws_BSM.Range("N79").Value = -17997.5
$ws_BSM.Range("H94").Value = 1562.75
$ws_BSM.Range("I94").Value = 1413.3043
$ws_BSM.Range("K94").Value = 1413.3043
$ws_BSM.Range("M94").Value = -962.3043
$ws_BSM.Range("H105").Value = 2752.1667
$ws_BSM.Range("I105").Value = 1254.5
$ws_BSM.Range("J105").Value = 5747.5
$ws_BSM.Range("K105").Value = 1254.5
$ws_BSM.Range("L105").Value = 5747.5
$ws_BSM.Range("M105").Value = 492.5
$ws_BSM.Range("N105").Value = -9241.5
$ws_BSM.Range("H134").Value = 9833.333000000001
$ws_BSM.Range("I134").Value = 14500
$ws_BSM.Range("J134").Value = 500
$ws_BSM.Range("K134").Value = 43500
$ws_BSM.Range("L134").Value = 1500
$ws_BSM.Range("M134").Value = -40965
$ws_BSM.Range("N134").Value = -6570

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H2").Value = 809.6
$ws_CRP.Range("I2").Value = 619.2
$ws_CRP.Range("J2").Value = 1000
$ws_CRP.Range("K2").Value = 619.2
$ws_CRP.Range("L2").Value = 1000
$ws_CRP.Range("M2").Value = -506.2
$ws_CRP.Range("N2").Value = -1226
$ws_CRP.Range("H22").Value = 900
$ws_CRP.Range("I22").Value = 900
$ws_CRP.Range("K22").Value = 900
$ws_CRP.Range("M22").Value = -550
$ws_CRP.Range("H31").Value = 6385.2856
$ws_CRP.Range("I31").Value = 5876.4
$ws_CRP.Range("K31").Value = 5876.4
$ws_CRP.Range("M31").Value = -5581.4
$ws_CRP.Range("H34").Value = 6385.2856
$ws_CRP.Range("I34").Value = 5876.4
$ws_CRP.Range("K34").Value = 5876.4
$ws_CRP.Range("M34").Value = -5674.4
$ws_CRP.Range("H134").Value = 10311.75
$ws_CRP.Range("I134").Value = 7125
$ws_CRP.Range("J134").Value = 13498.5
$ws_CRP.Range("K134").Value = 21375
$ws_CRP.Range("L134").Value = 40495.5
$ws_CRP.Range("M134").Value = -18840
$ws_CRP.Range("N134").Value = -45565.5

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H2").Value = 190.0625
$ws_CUL.Range("I2").Value = 194.45454
$ws_CUL.Range("J2").Value = 180.4
$ws_CUL.Range("K2").Value = 1166.72724
$ws_CUL.Range("L2").Value = 1082.4
$ws_CUL.Range("M2").Value = -1053.72724
$ws_CUL.Range("N2").Value = -1308.4
$ws_CUL.Range("H5").Value = 514.6316
$ws_CUL.Range("I5").Value = 525.7778
$ws_CUL.Range("K5").Value = 1577.3334
$ws_CUL.Range("M5").Value = -1465.3334
$ws_CUL.Range("H12").Value = 95.5
$ws_CUL.Range("I12").Value = 74
$ws_CUL.Range("J12").Value = 117
$ws_CUL.Range("K12").Value = 222
$ws_CUL.Range("L12").Value = 351
$ws_CUL.Range("M12").Value = -49
$ws_CUL.Range("N12").Value = -697
$ws_CUL.Range("H23").Value = 746.8
$ws_CUL.Range("I23").Value = 1149.5
$ws_CUL.Range("J23").Value = 646.125
$ws_CUL.Range("K23").Value = 3448.5
$ws_CUL.Range("L23").Value = 1938.375
$ws_CUL.Range("M23").Value = -3213.5
$ws_CUL.Range("N23").Value = -2408.375
$ws_CUL.Range("H37").Value = 65000
$ws_CUL.Range("J37").Value = 65000
$ws_CUL.Range("L37").Value = 195000
$ws_CUL.Range("N37").Value = -195224
$ws_CUL.Range("H40").Value = 150.26666
$ws_CUL.Range("J40").Value = 252
$ws_CUL.Range("L40").Value = 1008
$ws_CUL.Range("N40").Value = -1146
$ws_CUL.Range("H56").Value = 17854.615
$ws_CUL.Range("I56").Value = 17854.615
$ws_CUL.Range("K56").Value = 17854.615
$ws_CUL.Range("M56").Value = -17324.615
$ws_CUL.Range("H61").Value = 421.8
$ws_CUL.Range("J61").Value = 428
$ws_CUL.Range("L61").Value = 1284
$ws_CUL.Range("N61").Value = -1714
$ws_CUL.Range("H88").Value = 4914
$ws_CUL.Range("I88").Value = 4914
$ws_CUL.Range("K88").Value = 14742
$ws_CUL.Range("M88").Value = -14314
$ws_CUL.Range("H91").Value = 4914
$ws_CUL.Range("I91").Value = 4914
$ws_CUL.Range("K91").Value = 14742
$ws_CUL.Range("M91").Value = -13260
$ws_CUL.Range("H107").Value = 1380.25
$ws_CUL.Range("J107").Value = 1380.25
$ws_CUL.Range("L107").Value = 4140.75
$ws_CUL.Range("N107").Value = -7980.75
$ws_CUL.Range("H113").Value = 1699.5
$ws_CUL.Range("I113").Value = 999.5
$ws_CUL.Range("K113").Value = 2998.5
$ws_CUL.Range("M113").Value = -828.5
$ws_CUL.Range("H129").Value = 3976.3333
$ws_CUL.Range("J129").Value = 3976.3333
$ws_CUL.Range("L129").Value = 11928.9999
$ws_CUL.Range("N129").Value = -21928.9999
$ws_CUL.Range("H132").Value = 757.1667
$ws_CUL.Range("I132").Value = 757.1667
$ws_CUL.Range("K132").Value = 6814.5003
$ws_CUL.Range("M132").Value = -4284.5003
$ws_CUL.Range("H135").Value = 514.6316
$ws_CUL.Range("I135").Value = 525.7778
$ws_CUL.Range("K135").Value = 4732.000199999999
$ws_CUL.Range("M135").Value = -2197.000199999999
$ws_CUL.Range("H139").Value = 1258.8334
$ws_CUL.Range("I139").Value = 918.7273
$ws_CUL.Range("J139").Value = 5000
$ws_CUL.Range("K139").Value = 2756.1819
$ws_CUL.Range("L139").Value = 15000
$ws_CUL.Range("M139").Value = 2383.8181
$ws_CUL.Range("N139").Value = -25280

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H62").Value = 27000
$ws_GSM.Range("J62").Value = 0
$ws_GSM.Range("L62").Value = 0
$ws_GSM.Range("N62").ClearContents()
$ws_GSM.Range("H65").Value = 27000
$ws_GSM.Range("J65").Value = 0
$ws_GSM.Range("L65").Value = 0
$ws_GSM.Range("N65").ClearContents()
$ws_GSM.Range("H74").Value = 0
$ws_GSM.Range("J74").Value = 0
$ws_GSM.Range("L74").Value = 0
$ws_GSM.Range("N74").ClearContents()
$ws_GSM.Range("H77").Value = 0
$ws_GSM.Range("J77").Value = 0
$ws_GSM.Range("L77").Value = 0
$ws_GSM.Range("N77").ClearContents()
$ws_GSM.Range("H97").Value = 964.2222
$ws_GSM.Range("J97").Value = 1342
$ws_GSM.Range("L97").Value = 1342
$ws_GSM.Range("N97").Value = -2334
$ws_GSM.Range("H102").Value = 1728.1072
$ws_GSM.Range("I102").Value = 1803.7826
$ws_GSM.Range("J102").Value = 1380
$ws_GSM.Range("K102").Value = 1803.7826
$ws_GSM.Range("L102").Value = 1380
$ws_GSM.Range("M102").Value = -181.7826
$ws_GSM.Range("N102").Value = -4624
$ws_GSM.Range("H113").Value = 1133.2
$ws_GSM.Range("J113").Value = 0
$ws_GSM.Range("L113").Value = 0
$ws_GSM.Range("N113").ClearContents()

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H22").Value = 500
$ws_LTW.Range("J22").Value = 500
$ws_LTW.Range("L22").Value = 500
$ws_LTW.Range("N22").Value = -1090
$ws_LTW.Range("H27").Value = 500
$ws_LTW.Range("J27").Value = 500
$ws_LTW.Range("L27").Value = 500
$ws_LTW.Range("N27").Value = -714
$ws_LTW.Range("H46").Value = 3665
$ws_LTW.Range("I46").Value = 0
$ws_LTW.Range("J46").Value = 3665
$ws_LTW.Range("K46").Value = 0
$ws_LTW.Range("L46").Value = 3665
$ws_LTW.Range("M46").ClearContents()
$ws_LTW.Range("N46").Value = -4041

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H54").Value = 29833.334
$ws_WVR.Range("J54").Value = 29833.334
$ws_WVR.Range("L54").Value = 29833.334
$ws_WVR.Range("N54").Value = -30873.334
$ws_WVR.Range("H96").Value = 0
$ws_WVR.Range("I96").Value = 0
$ws_WVR.Range("K96").Value = 0
$ws_WVR.Range("M96").ClearContents()
